$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-11 follow the same visual pattern as row 5 (date / time-or-code / person / description / code / points).
# Copy the row-5 number formats down first so new cells land on the existing style indices
# (15 = date, 11 = time, 16 = general/text, 17 = general/points) instead of minting new styles.
for ($r = 6; $r -le 11; $r++) {
    $ws.Range("B5:G5").Copy() | Out-Null
    $ws.Range("B$r`:G$r").PasteSpecial(-4122) | Out-Null
}

# Rows 7-11 use the textual code "HD" in column C instead of a time-of-day fraction,
# so column C there should carry the same (General) format as column D, not the time format.
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C7:C11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Row 6 - Robert, logo animation
$ws.Range("B6").Value = 41964
$ws.Range("C6").Value = 0.375
$ws.Range("D6").Value = "Robert"
$ws.Range("E6").Value = "Animatie van het logo"
$ws.Range("F6").Value = 1001
$ws.Range("G6").Value = 13

# Row 7 - Robert, Jesse - sketches on paper (code "HD" entered after D/E,
# matching the shared-string insertion order in the original edit)
$ws.Range("B7").Value = 41964
$ws.Range("D7").Value = "Robert, Jesse"
$ws.Range("E7").Value = "Schetsen op papier van de game"
$ws.Range("F7").Value = 1028
$ws.Range("G7").Value = 40
$ws.Range("C7").Value = "HD"

# Row 8 - Jordy - use case diagrams
$ws.Range("B8").Value = 41964
$ws.Range("C8").Value = "HD"
$ws.Range("D8").Value = "Jordy"
$ws.Range("E8").Value = "Use case diagrammen"
$ws.Range("F8").Value = 1010
$ws.Range("G8").Value = 13

# Row 9 - Quincy - wireframes
$ws.Range("B9").Value = 41964
$ws.Range("C9").Value = "HD"
$ws.Range("D9").Value = "Quincy"
$ws.Range("E9").Value = "Wireframes"
$ws.Range("F9").Value = 1006
$ws.Range("G9").Value = 20

# Row 10 - Dean - HTML and CSS of homepage
$ws.Range("B10").Value = 41964
$ws.Range("C10").Value = "HD"
$ws.Range("D10").Value = "Dean"
$ws.Range("E10").Value = "HTML en CSS van de homepage"
$ws.Range("F10").Value = 1004
$ws.Range("G10").Value = 20

# Row 11 - Slawek - tagline and header
$ws.Range("B11").Value = 41964
$ws.Range("C11").Value = "HD"
$ws.Range("D11").Value = "Slawek"
$ws.Range("E11").Value = "verzinnen van de slagzin en de header"
$ws.Range("F11").Value = 1003
$ws.Range("G11").Value = 0.5

# Column E got wider to fit the longest description added above (~35.57 chars,
# matching Excel's "best fit" autosize for the new long description strings).
$ws.Columns("E").ColumnWidth = 34.65

# Selection moved to B12 after the new rows were entered.
$ws.Range("B12").Select()

